$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E retain text formatting so numeric-looking strings are not converted to numbers
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '42.367.15'
$ws.Range("E2").Value = '  -1.33%  '

# Row 3
$ws.Range("D3").Value = '2.299.27'
$ws.Range("E3").Value = '  -2.82%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '316.56'
$ws.Range("E5").Value = '  -0.05%  '

# Row 6
$ws.Range("D6").Value = '104.21'
$ws.Range("E6").Value = '  -4.85%  '

# Row 7
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  -1.09%  '

# Row 8
$ws.Range("E8").Value = '  +0.13%  '

# Row 9
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  -2.83%  '

# Row 10
$ws.Range("D10").Value = '39.72'
$ws.Range("E10").Value = '  -3.74%  '

# Row 11
$ws.Range("D11").Value = '0.0908'
$ws.Range("E11").Value = '  -2.35%  '

# Row 12
$ws.Range("D12").Value = '8.29'
$ws.Range("E12").Value = '  -4.10%  '

# Row 13
$ws.Range("E13").Value = '  -0.30%  '

# Row 14
$ws.Range("D14").Value = '0.964'
$ws.Range("E14").Value = '  -4.90%  '

# Row 15
$ws.Range("D15").Value = '15.33'
$ws.Range("E15").Value = '  -4.44%  '

# Row 16
$ws.Range("D16").Value = '2.646.85'
$ws.Range("E16").Value = '  -2.73%  '

# Row 17
$ws.Range("D17").Value = '2.293.38'
$ws.Range("E17").Value = '  -3.43%  '

# Row 18
$ws.Range("D18").Value = '42.259.28'
$ws.Range("E18").Value = '  -1.56%  '

# Row 19
$ws.Range("D19").Value = '7.41'
$ws.Range("E19").Value = '  -3.18%  '

# Row 20
$ws.Range("D20").Value = '0.0000106'
$ws.Range("E20").Value = '  -0.97%  '

# Row 21
$ws.Range("D21").Value = '73.41'
$ws.Range("E21").Value = '  -3.80%  '

# Row 22
$ws.Range("E22").Value = '  -0.77%  '

# Row 23
$ws.Range("D23").Value = '278.96'
$ws.Range("E23").Value = '  +4.37%  '

# Row 24
$ws.Range("D24").Value = '10.40'
$ws.Range("E24").Value = '  +9.41%  '

# Row 25
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  -2.49%  '

# Row 26
$ws.Range("E26").Value = '  +0.53%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '10.82'
$ws.Range("E27").Value = '  -5.79%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '2.37'
$ws.Range("E28").Value = '  +5.05%  '

# Row 29
$ws.Range("D29").Value = '22.86'
$ws.Range("E29").Value = '  -2.92%  '

# Row 30
$ws.Range("D30").Value = '36.23'
$ws.Range("E30").Value = '  -1.78%  '

# Row 31
$ws.Range("D31").Value = '164.14'
$ws.Range("E31").Value = '  -2.73%  '

# Row 32
$ws.Range("D32").Value = '0.0873'
$ws.Range("E32").Value = '  -4.28%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.82'
$ws.Range("E33").Value = '  -3.54%  '

# Row 34
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '2.83'
$ws.Range("E34").Value = '  -2.72%  '

# Row 35
$ws.Range("D35").Value = '0.136'
$ws.Range("E35").Value = '  +3.47%  '

# Row 36
$ws.Range("E36").Value = '  -4.53%  '

# Row 37
$ws.Range("D37").Value = '4.54'
$ws.Range("E37").Value = '  -4.45%  '

# Row 38
$ws.Range("D38").Value = '0.0348'
$ws.Range("E38").Value = '  -4.26%  '

# Row 39
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '3.76'
$ws.Range("E39").Value = '  -3.59%  '

# Row 40
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = '2.83'
$ws.Range("E40").Value = '  +4.38%  '

# Row 41
$ws.Range("D41").Value = '99.69'
$ws.Range("E41").Value = '  -5.90%  '

# Row 42
$ws.Range("D42").Value = '1.45'
$ws.Range("E42").Value = '  -4.28%  '

# Row 43
$ws.Range("D43").Value = '69.51'
$ws.Range("E43").Value = '  -2.60%  '

# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.226'
$ws.Range("E44").Value = '  -5.56%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.10%  '

# Row 46
$ws.Range("D46").Value = '12.01'
$ws.Range("E46").Value = '  -3.91%  '

# Row 47
$ws.Range("D47").Value = '112.58'
$ws.Range("E47").Value = '  -1.57%  '

# Row 48
$ws.Range("D48").Value = '77.78'
$ws.Range("E48").Value = '  -1.66%  '

# Row 49
$ws.Range("D49").Value = '8.92'
$ws.Range("E49").Value = '  -2.90%  '

# Row 50
$ws.Range("D50").Value = '5.30'
$ws.Range("E50").Value = '  -5.03%  '

# Row 51
$ws.Range("D51").Value = '1.587.93'
$ws.Range("E51").Value = '  +0.41%  '
